$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '1.000 (1.000 ± 0.000)'
$ws.Range("C2").Value = '00:00:44 (00:02:51 ± 00:01:34)'
$ws.Range("D2").Value = '00:00:05 (00:00:06 ± 00:00:00)'
$ws.Range("B3").Value = '1.000 (0.995 ± 0.005)'
$ws.Range("C3").Value = '00:00:07 (00:00:11 ± 00:00:03)'
$ws.Range("D3").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B4").Value = '1.000 (1.000 ± 0.000)'
$ws.Range("C4").Value = '00:00:30 (00:00:38 ± 00:00:06)'
$ws.Range("D4").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B5").Value = '1.000 (0.999 ± 0.003)'
$ws.Range("C5").Value = '00:05:06 (00:05:12 ± 00:00:03)'
$ws.Range("D5").Value = '00:00:01 (00:00:02 ± 00:00:00)'
$ws.Range("B6").Value = '1.000 (1.000 ± 0.001)'
$ws.Range("C6").Value = '00:04:57 (00:05:01 ± 00:00:02)'
$ws.Range("D6").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B7").Value = '0.956 (0.932 ± 0.013)'
$ws.Range("C7").Value = '00:05:00 (00:05:04 ± 00:00:02)'
$ws.Range("D7").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B8").Value = '1.000 (0.993 ± 0.010)'
$ws.Range("C8").Value = '00:02:47 (00:04:19 ± 00:01:02)'
$ws.Range("D8").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B9").Value = '1.000 (0.997 ± 0.004)'
$ws.Range("C9").Value = '00:04:59 (00:05:00 ± 00:00:00)'
$ws.Range("D9").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B10").Value = '1.000 (0.999 ± 0.002)'
$ws.Range("C10").Value = '00:04:29 (00:04:29 ± 00:00:00)'
$ws.Range("D10").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B11").Value = '0.890 (0.701 ± 0.154)'
$ws.Range("C11").Value = '00:05:01 (00:05:05 ± 00:00:00)'
$ws.Range("D11").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B12").Value = '0.463 (0.399 ± 0.030)'
$ws.Range("C12").Value = '00:02:53 (00:03:13 ± 00:00:07)'
$ws.Range("D12").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B13").Value = '1.000 (0.997 ± 0.005)'
$ws.Range("C13").Value = '00:00:02 (00:00:14 ± 00:00:07)'
$ws.Range("D13").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B14").Value = '1.000 (0.994 ± 0.006)'
$ws.Range("C14").Value = '00:00:20 (00:00:23 ± 00:00:02)'
$ws.Range("D14").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B15").Value = '1.000 (1.000 ± 0.000)'
$ws.Range("C15").Value = '00:00:00 (00:00:10 ± 00:00:06)'
$ws.Range("D15").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B16").Value = '1.000 (1.000 ± 0.001)'
$ws.Range("C16").Value = '00:00:08 (00:00:09 ± 00:00:00)'
$ws.Range("D16").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B17").Value = '1.000 (1.000 ± 0.000)'
$ws.Range("C17").Value = '00:05:01 (00:05:18 ± 00:00:13)'
$ws.Range("D17").Value = '00:00:00 (00:00:00 ± 00:00:00)'
